$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '328.85'
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '0.16%'
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '44.41'
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '1.07%'
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '5.601'
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.08075'
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '2.041'
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '7.17%'
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.9510'
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '0.94%'
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '2.518'
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '-8.50%'
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1166'
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '-1.20%'
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.1851'
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '-2.07%'
$cell.Style = "Normal"

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.09796'
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '2.82%'
$cell.Style = "Normal"

$ws.Range("B12").Value = 'BitrueCoin'

$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.04523'
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '7.94%'
$cell.Style = "Normal"

$ws.Range("B13").Value = 'BitMartToken'

$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.1069'
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '-0.14%'
$cell.Style = "Normal"

$ws.Range("B14").Value = 'BitForexToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.001277'
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '-0.51%'
$cell.Style = "Normal"

$ws.Range("B15").Value = 'CoinExToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.04202'
$cell.Style = "Normal"

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '-4.28%'
$cell.Style = "Normal"

$ws.Range("B16").Value = 'TigerCash'

$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.005904'
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '-0.75%'
$cell.Style = "Normal"

$ws.Range("B17").Value = 'LEO'

$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.368'
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '-5.25%'
$cell.Style = "Normal"

$ws.Range("B18").Value = 'GateToken'

$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '4.314'
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '0.22%'
$cell.Style = "Normal"

$ws.Range("B19").Value = 'BitpandaEcosystemToken'

$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.3482'
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '-0.99%'
$cell.Style = "Normal"

$ws.Range("B20").Value = 'MCDex'

$ws.Range("C20").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '10.24'
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '20.17%'
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.1411'
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '3.67%'
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '-3.91%'
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '0.43%'
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.004325'
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '0.33%'
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.0001191'
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '-3.83%'
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '-1.06%'
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02662'
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '-1.04%'
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.05551'
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '0.95%'
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.007575'
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '-2.84%'
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.1409'
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '1.13%'
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.008006'
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '-18.12%'
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.002018'
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '-5.09%'
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.008408'
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '-12.66%'
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.00007173'
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '0.78%'
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000752'
$cell.Style = "Normal"

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '-0.43%'
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.004007'
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '15.59%'
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.002275'
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '-0.59%'
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.00002105'
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '-0.43%'
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0002004'
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '-0.43%'
$cell.Style = "Normal"
